$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (shared-string reorder in the source diff is the
# visible side-effect of these two cells being retyped to the other
# country's name; the row data itself was separately refreshed below). ---
$ws.Range("A147").Value = "Burkina Faso"
$ws.Range("A148").Value = "Niger"

$ws.Range("A153").Value = "Malta"
$ws.Range("A154").Value = "Santo Tome y Principe"

# --- Last-updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 13:44"

# --- Updated case/recovered/death counts (columns B..H) ---
$ws.Range("E4").Value = 2275061
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 158375

$ws.Range("B21").Value = 211567
$ws.Range("C21").Value = 105
$ws.Range("E21").Value = 8741

$ws.Range("B27").Value = 111322
$ws.Range("C27").Value = 215
$ws.Range("D27").Value = 108002
$ws.Range("E27").Value = 3143

$ws.Range("B40").Value = 68299
$ws.Range("C40").Value = 388
$ws.Range("E40").Value = 8099

$ws.Range("B41").Value = 68166
$ws.Range("C41").Value = 99
$ws.Range("D41").Value = 62943
$ws.Range("E41").Value = 4652
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 571

$ws.Range("E52").Value = 2722
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 148

$ws.Range("B58").Value = 35616
$ws.Range("C58").Value = 66
$ws.Range("E58").Value = 2135

$ws.Range("B63").Value = 26066
$ws.Range("C63").Value = 730
$ws.Range("D63").Value = 16838
$ws.Range("E63").Value = 9071
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 157

$ws.Range("B67").Value = 21385
$ws.Range("C67").Value = 81
$ws.Range("D67").Value = 19063
$ws.Range("E67").Value = 1604

$ws.Range("B68").Value = 20750
$ws.Range("C68").Value = 418
$ws.Range("D68").Value = 14961
$ws.Range("E68").Value = 5732

$ws.Range("B73").Value = 17843
$ws.Range("C73").Value = 395
$ws.Range("D73").Value = 8809
$ws.Range("E73").Value = 8557

$ws.Range("B78").Value = 13996
$ws.Range("C78").Value = 207
$ws.Range("D78").Value = 12682
$ws.Range("E78").Value = 698
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 616

$ws.Range("B79").Value = 12541
$ws.Range("C79").Value = 244
$ws.Range("E79").Value = 7067

$ws.Range("B85").Value = 10386
$ws.Range("C85").Value = 42
$ws.Range("D85").Value = 6901
$ws.Range("E85").Value = 3274
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 211

$ws.Range("B96").Value = 6580
$ws.Range("C96").Value = 233
$ws.Range("D96").Value = 4701
$ws.Range("E96").Value = 1708
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 171

$ws.Range("B118").Value = 2824
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 2517
$ws.Range("E118").Value = 296

$ws.Range("D126").Value = 1831
$ws.Range("E126").Value = 228

$ws.Range("B147").Value = 1150
$ws.Range("C147").Value = 7
$ws.Range("D147").Value = 947
$ws.Range("E147").Value = 150
$ws.Range("H147").Value = 53

$ws.Range("B148").Value = 1147
$ws.Range("D148").Value = 1032
$ws.Range("E148").Value = 46
$ws.Range("H148").Value = 69

$ws.Range("C153").Value = 14
$ws.Range("D153").Value = 666
$ws.Range("E153").Value = 199
$ws.Range("H153").Value = 9

$ws.Range("B154").Value = 874
$ws.Range("D154").Value = 787
$ws.Range("E154").Value = 72
$ws.Range("H154").Value = 15

$ws.Range("B162").Value = 642
$ws.Range("C162").Value = 22
$ws.Range("E162").Value = 263
